# Add one more "Open issues" slide (Open issues (2)) at the end of the deck,
# using the same "Title and Content" layout (layout index 2) as the other
# body slides (e.g. the existing "Open issues" slide).

$p = $ppt.ActivePresentation

$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# Title
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Open issues (2)"

# Body content placeholder
$body = $s.Shapes.Item(2).TextFrame.TextRange

$body.Text = "Should controller software be able to name PSA_PORT_CPU and PSA_PORT_RECIRCULATE values in target-independent way?"

$p2 = $body.InsertAfter("`rSeems useful.  Seems like there should be a straightforward way to do this.")

$p3 = $body.InsertAfter("`rPotential quick hack:")

$p4a = $body.InsertAfter("`rTo generate P4Info file, compile with a psa.p4 file that has “")
$p4b = $body.InsertAfter("typedef")
$p4c = $body.InsertAfter(" bit<16> ")
$p4d = $body.InsertAfter("PortId_t")
$p4e = $body.InsertAfter(";”, the size in bits that you want ")
$p4f = $body.InsertAfter("PortId_t")
$p4g = $body.InsertAfter(" values to be between controllers/clients and agents/servers.")

$p5a = $body.InsertAfter("`rI do not think this is a good log term approach, because for the 2 PSA types with numerical translation, ")
$p5b = $body.InsertAfter("PortId_t")
$p5c = $body.InsertAfter(" and ")
$p5d = $body.InsertAfter("ClassOfService_t")
$p5e = $body.InsertAfter(", we still want to minimize annotations in P4 program myp4prog.p4 wherever those types are used")
$p5f = $body.InsertAfter(", ideally 0 of them.")

# Outline levels: paragraphs 2, 4, 5 are at level 2 (index 1 / lvl="1")
$body.Paragraphs(2,1).IndentLevel = 2
$body.Paragraphs(4,1).IndentLevel = 2
$body.Paragraphs(5,1).IndentLevel = 2
